$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 948.1111
$ws.Range("I15").Value = 948.1111
$ws.Range("K15").Value = 2844.3333
$ws.Range("M15").Value = -2675.3333
$ws.Range("H17").Value = 989.46
$ws.Range("J17").Value = 1005.93475
$ws.Range("L17").Value = 3017.80425
$ws.Range("N17").Value = -3353.80425
$ws.Range("H19").Value = 2437.5715
$ws.Range("I19").Value = 2314.8
$ws.Range("K19").Value = 2314.8
$ws.Range("M19").Value = -2139.8
$ws.Range("H33").Value = 583.8333
$ws.Range("I33").Value = 289
$ws.Range("J33").Value = 878.6667
$ws.Range("K33").Value = 289
$ws.Range("L33").Value = 878.6667
$ws.Range("M33").Value = -60
$ws.Range("N33").Value = -1336.6667
$ws.Range("H43").Value = 3368.4
$ws.Range("I43").Value = 3111.75
$ws.Range("K43").Value = 3111.75
$ws.Range("M43").Value = -3042.75
$ws.Range("H57").Value = 68840
$ws.Range("J57").Value = 68840
$ws.Range("L57").Value = 206520
$ws.Range("N57").Value = -207518
$ws.Range("H70").Value = 1043.25
$ws.Range("I70").Value = 799.3570999999999
$ws.Range("J70").Value = 2750.5
$ws.Range("K70").Value = 2398.0713
$ws.Range("L70").Value = 8251.5
$ws.Range("M70").Value = -2128.0713
$ws.Range("N70").Value = -8791.5
$ws.Range("H73").Value = 1043.25
$ws.Range("I73").Value = 799.3570999999999
$ws.Range("J73").Value = 2750.5
$ws.Range("K73").Value = 2398.0713
$ws.Range("L73").Value = 8251.5
$ws.Range("M73").Value = -1462.0713
$ws.Range("N73").Value = -10123.5
$ws.Range("H74").Value = 3938
$ws.Range("I74").Value = 3245.0908
$ws.Range("K74").Value = 3245.0908
$ws.Range("M74").Value = -2309.0908
$ws.Range("H77").Value = 3938
$ws.Range("I77").Value = 3245.0908
$ws.Range("K77").Value = 16225.454
$ws.Range("M77").Value = -11545.454
$ws.Range("H80").Value = 3920.6
$ws.Range("I80").Value = 523.8
$ws.Range("K80").Value = 1571.4
$ws.Range("M80").Value = -573.3999999999999
$ws.Range("H83").Value = 3920.6
$ws.Range("I83").Value = 523.8
$ws.Range("K83").Value = 4714.2
$ws.Range("M83").Value = 277.8000000000002
$ws.Range("H86").Value = 679.625
$ws.Range("I86").Value = 680
$ws.Range("J86").Value = 679
$ws.Range("K86").Value = 680
$ws.Range("L86").Value = 679
$ws.Range("M86").Value = 443
$ws.Range("N86").Value = -2925
$ws.Range("H89").Value = 679.625
$ws.Range("I89").Value = 680
$ws.Range("J89").Value = 679
$ws.Range("K89").Value = 3400
$ws.Range("L89").Value = 3395
$ws.Range("M89").Value = 2216
$ws.Range("N89").Value = -14627
$ws.Range("H96").Value = 341.85715
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").Value = ""
$ws.Range("H99").Value = 1603.9286
$ws.Range("I99").Value = 597.2222
$ws.Range("J99").Value = 3416
$ws.Range("K99").Value = 1791.6666
$ws.Range("L99").Value = 10248
$ws.Range("M99").Value = -293.6666
$ws.Range("N99").Value = -13244
$ws.Range("H101").Value = 2377.1667
$ws.Range("I101").Value = 688
$ws.Range("J101").Value = 4066.3333
$ws.Range("K101").Value = 2064
$ws.Range("L101").Value = 12198.9999
$ws.Range("M101").Value = -442
$ws.Range("N101").Value = -15442.9999
$ws.Range("H107").Value = 370.33334
$ws.Range("I107").Value = 369
$ws.Range("K107").Value = 369
$ws.Range("M107").Value = 1551
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = ""
$ws.Range("H112").Value = 1450.1875
$ws.Range("J112").Value = 1450.1875
$ws.Range("L112").Value = 4350.5625
$ws.Range("N112").Value = -6566.5625
$ws.Range("H113").Value = 2253.125
$ws.Range("I113").Value = 2304.2856
$ws.Range("J113").Value = 1895
$ws.Range("K113").Value = 2304.2856
$ws.Range("L113").Value = 1895
$ws.Range("M113").Value = 949.7143999999998
$ws.Range("N113").Value = -8403
$ws.Range("H114").Value = 70189
$ws.Range("I114").Value = 69656
$ws.Range("J114").Value = 70722
$ws.Range("K114").Value = 69656
$ws.Range("L114").Value = 70722
$ws.Range("M114").Value = -65317
$ws.Range("N114").Value = -79400
$ws.Range("H115").Value = 472.85715
$ws.Range("I115").Value = 472.85715
$ws.Range("K115").Value = 1418.57145
$ws.Range("M115").Value = 148.4285500000001
$ws.Range("H116").Value = 99270.55
$ws.Range("I116").Value = 9327.333000000001
$ws.Range("K116").Value = 9327.333000000001
$ws.Range("M116").Value = -5885.333000000001
$ws.Range("H125").Value = 860.8333
$ws.Range("I125").Value = 803.3333
$ws.Range("J125").Value = 880
$ws.Range("K125").Value = 7229.9997
$ws.Range("L125").Value = 7920
$ws.Range("M125").Value = -4769.9997
$ws.Range("N125").Value = -12840
$ws.Range("H132").Value = 2160.7222
$ws.Range("I132").Value = 2035.9286
$ws.Range("K132").Value = 6107.7858
$ws.Range("M132").Value = -3577.7858
$ws.Range("H135").Value = 1014.2273
$ws.Range("I135").Value = 606.82355
$ws.Range("K135").Value = 5461.41195
$ws.Range("M135").Value = -2926.41195
$ws.Range("H137").Value = 2498.8408
$ws.Range("I137").Value = 1776.2646
$ws.Range("K137").Value = 5328.793799999999
$ws.Range("M137").Value = -2778.793799999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3970.6
$ws.Range("J2").Value = 11513
$ws.Range("L2").Value = 11513
$ws.Range("N2").Value = -11739
$ws.Range("H5").Value = 121
$ws.Range("I5").Value = 104.8
$ws.Range("J5").Value = 202
$ws.Range("K5").Value = 104.8
$ws.Range("L5").Value = 202
$ws.Range("M5").Value = 7.200000000000003
$ws.Range("N5").Value = -426
$ws.Range("H32").Value = 1710.5428
$ws.Range("I32").Value = 1710.5428
$ws.Range("K32").Value = 1710.5428
$ws.Range("M32").Value = -1423.5428
$ws.Range("H61").Value = 2440.5
$ws.Range("I61").Value = 2134.6
$ws.Range("K61").Value = 2134.6
$ws.Range("M61").Value = -1922.6
$ws.Range("H74").Value = 2725179.8
$ws.Range("I74").Value = 1426358.9
$ws.Range("J74").Value = 6946347.5
$ws.Range("K74").Value = 1426358.9
$ws.Range("L74").Value = 6946347.5
$ws.Range("M74").Value = -1425484.9
$ws.Range("N74").Value = -6948095.5
$ws.Range("H77").Value = 2725179.8
$ws.Range("I77").Value = 1426358.9
$ws.Range("J77").Value = 6946347.5
$ws.Range("K77").Value = 7131794.5
$ws.Range("L77").Value = 34731737.5
$ws.Range("M77").Value = -7127426.5
$ws.Range("N77").Value = -34740473.5
$ws.Range("H97").Value = 674.55554
$ws.Range("I97").Value = 469.2857
$ws.Range("K97").Value = 469.2857
$ws.Range("M97").Value = 26.71429999999998
$ws.Range("H110").Value = 2306.2122
$ws.Range("I110").Value = 1799.84
$ws.Range("K110").Value = 1799.84
$ws.Range("M110").Value = 245.1600000000001
$ws.Range("H116").Value = 3970.6
$ws.Range("J116").Value = 11513
$ws.Range("L116").Value = 11513
$ws.Range("N116").Value = -16101
$ws.Range("H122").Value = 326.2
$ws.Range("I122").Value = 326.2
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 978.5999999999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 1471.4
$ws.Range("N122").Value = ""
$ws.Range("H132").Value = 31254870
$ws.Range("I132").Value = 4194.12
$ws.Range("J132").Value = 142864430
$ws.Range("K132").Value = 12582.36
$ws.Range("L132").Value = 428593290
$ws.Range("M132").Value = -10052.36
$ws.Range("N132").Value = -428598350
$ws.Range("H136").Value = 2440.5
$ws.Range("I136").Value = 2134.6
$ws.Range("K136").Value = 6403.799999999999
$ws.Range("M136").Value = -3853.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3970.6
$ws.Range("J3").Value = 11513
$ws.Range("L3").Value = 11513
$ws.Range("N3").Value = -11741
$ws.Range("H4").Value = 121
$ws.Range("I4").Value = 104.8
$ws.Range("J4").Value = 202
$ws.Range("K4").Value = 104.8
$ws.Range("L4").Value = 202
$ws.Range("M4").Value = 10.2
$ws.Range("N4").Value = -432
$ws.Range("H5").Value = 598.6
$ws.Range("I5").Value = 1306.5
$ws.Range("J5").Value = 126.666664
$ws.Range("K5").Value = 1306.5
$ws.Range("L5").Value = 126.666664
$ws.Range("M5").Value = -1193.5
$ws.Range("N5").Value = -352.666664
$ws.Range("H20").Value = 3121.9443
$ws.Range("I20").Value = 2308.4167
$ws.Range("J20").Value = 4749
$ws.Range("K20").Value = 2308.4167
$ws.Range("L20").Value = 4749
$ws.Range("M20").Value = -2061.4167
$ws.Range("N20").Value = -5243
$ws.Range("H86").Value = 3022.6667
$ws.Range("I86").Value = 1991.875
$ws.Range("J86").Value = 4200.7144
$ws.Range("K86").Value = 1991.875
$ws.Range("L86").Value = 4200.7144
$ws.Range("M86").Value = -868.875
$ws.Range("N86").Value = -6446.7144
$ws.Range("H88").Value = 28447.666
$ws.Range("J88").Value = 28447.666
$ws.Range("L88").Value = 28447.666
$ws.Range("N88").Value = -29259.666
$ws.Range("H89").Value = 3022.6667
$ws.Range("I89").Value = 1991.875
$ws.Range("J89").Value = 4200.7144
$ws.Range("K89").Value = 9959.375
$ws.Range("L89").Value = 21003.572
$ws.Range("M89").Value = -4343.375
$ws.Range("N89").Value = -32235.572
$ws.Range("H91").Value = 28447.666
$ws.Range("J91").Value = 28447.666
$ws.Range("L91").Value = 28447.666
$ws.Range("N91").Value = -31255.666
$ws.Range("H92").Value = 32566.334
$ws.Range("J92").Value = 32566.334
$ws.Range("L92").Value = 32566.334
$ws.Range("N92").Value = -37558.334
$ws.Range("H94").Value = 15794.462
$ws.Range("I94").Value = 22604.334
$ws.Range("J94").Value = 472.25
$ws.Range("K94").Value = 22604.334
$ws.Range("L94").Value = 472.25
$ws.Range("M94").Value = -22153.334
$ws.Range("N94").Value = -1374.25
$ws.Range("H99").Value = 4917.25
$ws.Range("I99").Value = 4943.5713
$ws.Range("J99").Value = 4733
$ws.Range("K99").Value = 4943.5713
$ws.Range("L99").Value = 4733
$ws.Range("M99").Value = -3445.5713
$ws.Range("N99").Value = -7729
$ws.Range("H105").Value = 3541.2104
$ws.Range("I105").Value = 3457.1765
$ws.Range("K105").Value = 3457.1765
$ws.Range("M105").Value = -1710.1765
$ws.Range("H134").Value = 9264115
$ws.Range("I134").Value = 10005026
$ws.Range("K134").Value = 30015078
$ws.Range("M134").Value = -30012543

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 3566.6667
$ws.Range("I10").Value = 5250
$ws.Range("J10").Value = 200
$ws.Range("K10").Value = 5250
$ws.Range("L10").Value = 200
$ws.Range("M10").Value = -5111
$ws.Range("N10").Value = -478
$ws.Range("H13").Value = 7201
$ws.Range("I13").Value = 1899
$ws.Range("J13").Value = 8968.333000000001
$ws.Range("K13").Value = 1899
$ws.Range("L13").Value = 8968.333000000001
$ws.Range("M13").Value = -1760
$ws.Range("N13").Value = -9246.333000000001
$ws.Range("H16").Value = 1947.8
$ws.Range("I16").Value = 1970.75
$ws.Range("K16").Value = 1970.75
$ws.Range("M16").Value = -1683.75
$ws.Range("H31").Value = 2688.647
$ws.Range("I31").Value = 902.2727
$ws.Range("J31").Value = 5963.6665
$ws.Range("K31").Value = 902.2727
$ws.Range("L31").Value = 5963.6665
$ws.Range("M31").Value = -607.2727
$ws.Range("N31").Value = -6553.6665
$ws.Range("H34").Value = 2688.647
$ws.Range("I34").Value = 902.2727
$ws.Range("J34").Value = 5963.6665
$ws.Range("K34").Value = 902.2727
$ws.Range("L34").Value = 5963.6665
$ws.Range("M34").Value = -700.2727
$ws.Range("N34").Value = -6367.6665
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").Value = ""
$ws.Range("H58").Value = 3121.3333
$ws.Range("I58").Value = 3121.3333
$ws.Range("K58").Value = 3121.3333
$ws.Range("M58").Value = -2918.3333
$ws.Range("H62").Value = 20967.334
$ws.Range("I62").Value = 17789
$ws.Range("J62").Value = 24145.666
$ws.Range("K62").Value = 17789
$ws.Range("L62").Value = 24145.666
$ws.Range("M62").Value = -17165
$ws.Range("N62").Value = -25393.666
$ws.Range("H65").Value = 20967.334
$ws.Range("I65").Value = 17789
$ws.Range("J65").Value = 24145.666
$ws.Range("K65").Value = 88945
$ws.Range("L65").Value = 120728.33
$ws.Range("M65").Value = -85825
$ws.Range("N65").Value = -126968.33
$ws.Range("H99").Value = 4817.8184
$ws.Range("I99").Value = 4817.8184
$ws.Range("K99").Value = 4817.8184
$ws.Range("M99").Value = -3319.8184
$ws.Range("H107").Value = 2152.2273
$ws.Range("I107").Value = 1385.5
$ws.Range("J107").Value = 3072.3
$ws.Range("K107").Value = 1385.5
$ws.Range("L107").Value = 3072.3
$ws.Range("M107").Value = 534.5
$ws.Range("N107").Value = -6912.3
$ws.Range("H112").Value = 53000
$ws.Range("J112").Value = 53000
$ws.Range("L112").Value = 53000
$ws.Range("N112").Value = -55954
$ws.Range("H113").Value = 1947.8
$ws.Range("I113").Value = 1970.75
$ws.Range("K113").Value = 1970.75
$ws.Range("M113").Value = 199.25
$ws.Range("H122").Value = 2299.2856
$ws.Range("I122").Value = 2250.25
$ws.Range("J122").Value = 2364.6667
$ws.Range("K122").Value = 6750.75
$ws.Range("L122").Value = 7094.000100000001
$ws.Range("M122").Value = -4300.75
$ws.Range("N122").Value = -11994.0001
$ws.Range("H126").Value = 4817.8184
$ws.Range("I126").Value = 4817.8184
$ws.Range("K126").Value = 14453.4552
$ws.Range("M126").Value = -11983.4552
$ws.Range("H134").Value = 3590.158
$ws.Range("J134").Value = 5648.1665
$ws.Range("L134").Value = 16944.4995
$ws.Range("N134").Value = -22014.4995
$ws.Range("H136").Value = 3121.3333
$ws.Range("I136").Value = 3121.3333
$ws.Range("K136").Value = 9363.999899999999
$ws.Range("M136").Value = -6813.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1175
$ws.Range("J2").Value = 1175
$ws.Range("L2").Value = 7050
$ws.Range("N2").Value = -7276
$ws.Range("H4").Value = 131126.27
$ws.Range("J4").Value = 334944.88
$ws.Range("L4").Value = 1004834.64
$ws.Range("N4").Value = -1005058.64
$ws.Range("H7").Value = 627.4167
$ws.Range("I7").Value = 147.44444
$ws.Range("J7").Value = 2067.3333
$ws.Range("K7").Value = 442.33332
$ws.Range("L7").Value = 6201.999899999999
$ws.Range("M7").Value = -330.33332
$ws.Range("N7").Value = -6425.999899999999
$ws.Range("H16").Value = 200
$ws.Range("I16").Value = 200
$ws.Range("K16").Value = 600
$ws.Range("M16").Value = -427
$ws.Range("H23").Value = 1481.8
$ws.Range("I23").Value = 2408
$ws.Range("J23").Value = 555.6
$ws.Range("K23").Value = 7224
$ws.Range("L23").Value = 1666.8
$ws.Range("M23").Value = -6989
$ws.Range("N23").Value = -2136.8
$ws.Range("H34").Value = 4789
$ws.Range("I34").Value = 173.66667
$ws.Range("K34").Value = 521.00001
$ws.Range("M34").Value = -437.00001
$ws.Range("H39").Value = 7063.7144
$ws.Range("J39").Value = 7260.923
$ws.Range("L39").Value = 21782.769
$ws.Range("N39").Value = -22370.769
$ws.Range("H54").Value = 17295.666
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").Value = ""
$ws.Range("H107").Value = 499.26666
$ws.Range("J107").Value = 630.25
$ws.Range("L107").Value = 1890.75
$ws.Range("N107").Value = -5730.75
$ws.Range("H109").Value = 4966
$ws.Range("I109").Value = 4949
$ws.Range("K109").Value = 14847
$ws.Range("M109").Value = -13807
$ws.Range("H110").Value = 3000
$ws.Range("I110").Value = 3000
$ws.Range("K110").Value = 9000
$ws.Range("M110").Value = -4910
$ws.Range("H112").Value = 2994.3333
$ws.Range("I112").Value = 449.75
$ws.Range("K112").Value = 1349.25
$ws.Range("M112").Value = -241.25
$ws.Range("H114").Value = 2147.1538
$ws.Range("J114").Value = 3131
$ws.Range("L114").Value = 9393
$ws.Range("N114").Value = -15901
$ws.Range("H115").Value = 2968
$ws.Range("I115").Value = 2452.25
$ws.Range("J115").Value = 3999.5
$ws.Range("K115").Value = 7356.75
$ws.Range("L115").Value = 11998.5
$ws.Range("M115").Value = -6181.75
$ws.Range("N115").Value = -14348.5
$ws.Range("H116").Value = 43100.64
$ws.Range("I116").Value = 85207.35000000001
$ws.Range("J116").Value = 10563.637
$ws.Range("K116").Value = 255622.05
$ws.Range("L116").Value = 31690.911
$ws.Range("M116").Value = -252180.05
$ws.Range("N116").Value = -38574.911
$ws.Range("H118").Value = 2000
$ws.Range("J118").Value = 2000
$ws.Range("L118").Value = 6000
$ws.Range("N118").Value = -8486
$ws.Range("H119").Value = 3498.3333
$ws.Range("I119").Value = 3498.3333
$ws.Range("K119").Value = 10494.9999
$ws.Range("M119").Value = -5656.999899999999
$ws.Range("H128").Value = 549895
$ws.Range("I128").Value = 549895
$ws.Range("K128").Value = 1649685
$ws.Range("M128").Value = -1644705
$ws.Range("H131").Value = 368190.72
$ws.Range("J131").Value = 469280.9
$ws.Range("L131").Value = 1407842.7
$ws.Range("N131").Value = -1417922.7
$ws.Range("H132").Value = 1221.4286
$ws.Range("I132").Value = 301.5
$ws.Range("J132").Value = 1589.4
$ws.Range("K132").Value = 2713.5
$ws.Range("L132").Value = 14304.6
$ws.Range("M132").Value = -183.5
$ws.Range("N132").Value = -19364.6
$ws.Range("H137").Value = 398087.22
$ws.Range("J137").Value = 580096.75
$ws.Range("L137").Value = 1740290.25
$ws.Range("N137").Value = -1750490.25
$ws.Range("H140").Value = 3935.8462
$ws.Range("I140").Value = 3966.6
$ws.Range("K140").Value = 11899.8
$ws.Range("M140").Value = -6719.799999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 22559.6
$ws.Range("I59").Value = 24449
$ws.Range("J59").Value = 21300
$ws.Range("K59").Value = 24449
$ws.Range("L59").Value = 21300
$ws.Range("M59").Value = -23866
$ws.Range("N59").Value = -22466
$ws.Range("H70").Value = 7257.567
$ws.Range("I70").Value = 7337.857
$ws.Range("K70").Value = 7337.857
$ws.Range("M70").Value = -7067.857
$ws.Range("H73").Value = 7257.567
$ws.Range("I73").Value = 7337.857
$ws.Range("K73").Value = 7337.857
$ws.Range("M73").Value = -6401.857
$ws.Range("H80").Value = 2401.4666
$ws.Range("I80").Value = 2367.3
$ws.Range("K80").Value = 2367.3
$ws.Range("M80").Value = -1369.3
$ws.Range("H83").Value = 2401.4666
$ws.Range("I83").Value = 2367.3
$ws.Range("K83").Value = 11836.5
$ws.Range("M83").Value = -6844.5
$ws.Range("H97").Value = 537.1053000000001
$ws.Range("I97").Value = 580.2143
$ws.Range("J97").Value = 416.4
$ws.Range("K97").Value = 580.2143
$ws.Range("L97").Value = 416.4
$ws.Range("M97").Value = -84.21429999999998
$ws.Range("N97").Value = -1408.4
$ws.Range("H102").Value = 2446.9048
$ws.Range("I102").Value = 2206.3333
$ws.Range("J102").Value = 2543.1333
$ws.Range("K102").Value = 2206.3333
$ws.Range("L102").Value = 2543.1333
$ws.Range("M102").Value = -584.3332999999998
$ws.Range("N102").Value = -5787.1333
$ws.Range("H113").Value = 2930.6
$ws.Range("I113").Value = 2881.4
$ws.Range("K113").Value = 2881.4
$ws.Range("M113").Value = -711.4000000000001
$ws.Range("H122").Value = 2240.9534
$ws.Range("I122").Value = 2212.3125
$ws.Range("K122").Value = 6636.9375
$ws.Range("M122").Value = -4186.9375
$ws.Range("H126").Value = 5467.7
$ws.Range("I126").Value = 16199.5
$ws.Range("J126").Value = 2784.75
$ws.Range("K126").Value = 48598.5
$ws.Range("L126").Value = 8354.25
$ws.Range("M126").Value = -46128.5
$ws.Range("N126").Value = -13294.25
$ws.Range("H132").Value = 1280.4286
$ws.Range("I132").Value = 1286.5555
$ws.Range("J132").Value = 1269.4
$ws.Range("K132").Value = 3859.6665
$ws.Range("L132").Value = 3808.2
$ws.Range("M132").Value = -1329.6665
$ws.Range("N132").Value = -8868.200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3894.2
$ws.Range("J7").Value = 3756.7778
$ws.Range("L7").Value = 3756.7778
$ws.Range("N7").Value = -3980.7778
$ws.Range("H40").Value = 2946.25
$ws.Range("I40").Value = 2352.1428
$ws.Range("J40").Value = 3778
$ws.Range("K40").Value = 2352.1428
$ws.Range("L40").Value = 3778
$ws.Range("M40").Value = -2216.1428
$ws.Range("N40").Value = -4050
$ws.Range("H61").Value = 6398.8
$ws.Range("I61").Value = 5331.6665
$ws.Range("K61").Value = 5331.6665
$ws.Range("M61").Value = -5129.6665
$ws.Range("H82").Value = 6976.3076
$ws.Range("I82").Value = 4116.5
$ws.Range("J82").Value = 9427.571
$ws.Range("K82").Value = 4116.5
$ws.Range("L82").Value = 9427.571
$ws.Range("M82").Value = -3755.5
$ws.Range("N82").Value = -10149.571
$ws.Range("H85").Value = 6976.3076
$ws.Range("I85").Value = 4116.5
$ws.Range("J85").Value = 9427.571
$ws.Range("K85").Value = 4116.5
$ws.Range("L85").Value = 9427.571
$ws.Range("M85").Value = -2868.5
$ws.Range("N85").Value = -11923.571
$ws.Range("H93").Value = 1345.8695
$ws.Range("I93").Value = 830.6667
$ws.Range("K93").Value = 830.6667
$ws.Range("M93").Value = 417.3333
$ws.Range("H113").Value = 6398.8
$ws.Range("I113").Value = 5331.6665
$ws.Range("K113").Value = 5331.6665
$ws.Range("M113").Value = -3161.6665
$ws.Range("H122").Value = 3458.6538
$ws.Range("I122").Value = 3163.8333
$ws.Range("K122").Value = 9491.499899999999
$ws.Range("M122").Value = -7041.499899999999
$ws.Range("H126").Value = 3894.2
$ws.Range("J126").Value = 3756.7778
$ws.Range("L126").Value = 11270.3334
$ws.Range("N126").Value = -16210.3334
$ws.Range("H132").Value = 2733.1177
$ws.Range("I132").Value = 2779.25
$ws.Range("J132").Value = 1995
$ws.Range("K132").Value = 8337.75
$ws.Range("L132").Value = 5985
$ws.Range("M132").Value = -5807.75
$ws.Range("N132").Value = -11045

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 24998.666
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 24998.666
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 24998.666
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = -25222.666
$ws.Range("H4").Value = 3835.4546
$ws.Range("J4").Value = 3835.4546
$ws.Range("L4").Value = 3835.4546
$ws.Range("N4").Value = -4061.4546
$ws.Range("H46").Value = 38499
$ws.Range("I46").Value = 47998
$ws.Range("J46").Value = 29000
$ws.Range("K46").Value = 47998
$ws.Range("L46").Value = 29000
$ws.Range("M46").Value = -47767
$ws.Range("N46").Value = -29462
$ws.Range("H51").Value = 40077
$ws.Range("J51").Value = 40077
$ws.Range("L51").Value = 40077
$ws.Range("N51").Value = -41097
$ws.Range("H62").Value = 3520.2222
$ws.Range("I62").Value = 3345
$ws.Range("K62").Value = 3345
$ws.Range("M62").Value = -2721
$ws.Range("H65").Value = 3520.2222
$ws.Range("I65").Value = 3345
$ws.Range("K65").Value = 16725
$ws.Range("M65").Value = -13605
$ws.Range("H81").Value = 10204.412
$ws.Range("I81").Value = 12444
$ws.Range("J81").Value = 6098.5
$ws.Range("K81").Value = 24888
$ws.Range("L81").Value = 12197
$ws.Range("M81").Value = -23827
$ws.Range("N81").Value = -14319
$ws.Range("H84").Value = 10204.412
$ws.Range("I84").Value = 12444
$ws.Range("J84").Value = 6098.5
$ws.Range("K84").Value = 124440
$ws.Range("L84").Value = 60985
$ws.Range("M84").Value = -119136
$ws.Range("N84").Value = -71593
$ws.Range("H111").Value = 40000
$ws.Range("J111").Value = 40000
$ws.Range("L111").Value = 40000
$ws.Range("N111").Value = -48180
$ws.Range("H113").Value = 3596.4285
$ws.Range("I113").Value = 919.4286
$ws.Range("J113").Value = 6273.4287
$ws.Range("K113").Value = 2758.2858
$ws.Range("L113").Value = 18820.2861
$ws.Range("M113").Value = -588.2857999999997
$ws.Range("N113").Value = -23160.2861
$ws.Range("H122").Value = 2321.5278
$ws.Range("I122").Value = 2793.4211
$ws.Range("J122").Value = 1794.1177
$ws.Range("K122").Value = 8380.263300000001
$ws.Range("L122").Value = 5382.3531
$ws.Range("M122").Value = -5930.263300000001
$ws.Range("N122").Value = -10282.3531
$ws.Range("H132").Value = 1744.4286
$ws.Range("I132").Value = 1664.9474
$ws.Range("J132").Value = 2499.5
$ws.Range("K132").Value = 4994.8422
$ws.Range("L132").Value = 7498.5
$ws.Range("M132").Value = -2464.8422
$ws.Range("N132").Value = -12558.5
$ws.Range("H134").Value = 38499
$ws.Range("I134").Value = 47998
$ws.Range("J134").Value = 29000
$ws.Range("K134").Value = 143994
$ws.Range("L134").Value = 87000
$ws.Range("M134").Value = -141459
$ws.Range("N134").Value = -92070
